$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rnaDate: "6.05.19" -> "06.05.19" (rows 2-27, column A) ---
# Force text format while typing so Excel doesn't re-parse the dotted
# date-looking string back into a date serial, then restore General so the
# cell's number format / style matches the rest of the column.
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "06.05.19"
    $cell.NumberFormat = "General"
}

# --- polyAIsolationProtocol: kit code relabel, all rows now "...E7490L" ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 7).Value = "NEBNextPoly(A)E7490L"
}

# --- roboticS1Prep: plain boolean literal -> explicit =FALSE() formula ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=FALSE()"
}

# --- column widths: G/H/I get their own explicit widths (target stored
# widths of 21.91 / 10.61 / 13.42 chars; Excel's COM ColumnWidth snaps to a
# pixel grid so these are the closest reachable values) ---
$ws.Columns.Item(7).ColumnWidth = 20.92
$ws.Columns.Item(8).ColumnWidth = 9.76
$ws.Columns.Item(9).ColumnWidth = 12.6

# --- selection moves from I27 to A11 ---
$ws.Range("A11").Select()
